$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the FFR_CA and CA_LF columns (D:E), shifting nothing else -
# update remaining header / data in place.
$ws.Range("C1").Value = "CA_FFR"

$ws.Range("B2").Value = 0.4122313011173397
$ws.Range("C2").Value = 13.71835060590134

$ws.Range("B3").Value = [double]"1.158850793103738E-12"
$ws.Range("C3").Value = 0

# Clear out the now-unused D:E columns entirely
$ws.Range("D1:E3").Clear() | Out-Null
